$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item(1)
$wsZhCn = $wb.Worksheets.Item(2)
$wsDeDe = $wb.Worksheets.Item(3)

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e186cd9da09cb534ea01c7b65d1224661e356e4/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2e186cd9da09cb534ea01c7b65d1224661e356e4/e2e/b.md"

$newStatus = "Handed back: in sync with en-US"

# 1. Update status text everywhere it is shown (Overview summary columns + per-language Status column)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# 2. Fill in handback report info for zh-cn sheet
$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-24 14:40:58"

$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-24 14:40:58"

# 3. Fill in handback report info for de-de sheet
$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-24 14:41:18"

$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-24 14:41:18"

# 4. Rebuild hyperlinks on zh-cn and de-de sheets so that the new "Latest Target File" (I) cells
#    become hyperlinks too, in the order: A2, I2, A3, I3
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlA, $null, $null, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlA, $null, $null, "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlB, $null, $null, "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlA, $null, $null, "a.md")

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlA, $null, $null, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlA, $null, $null, "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlB, $null, $null, "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlA, $null, $null, "a.md")

# 5. Column width adjustments (widen Status columns and Latest Handback File columns)
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1
